# DemoQA test workbook: keyword-driven framework tweak.
#
# - "Test Steps" sheet, column H ("Data") rows 2 & 3 previously held the
#   per-step keyword values "browser" / "url" — these are now just "NA"
#   (that data moved / isn't used that way any more).
# - Row 12's Data cell picks up what used to be the "fileLocation" value,
#   now reusing the "browser" keyword text instead.
# - "Test Steps" becomes the active/selected sheet (instead of "Test
#   Cases"), with the selection parked on H1.

$wb = $excel.ActiveWorkbook

$testSteps = $wb.Worksheets.Item("Test Steps")

$testSteps.Range("H2").Value = "NA"
$testSteps.Range("H3").Value = "NA"
$testSteps.Range("H12").Value = "browser"

$testSteps.Activate()
$testSteps.Range("H1").Select()
